# The "reviews_count" column (column E) was removed from the sheet,
# so every column to its right (reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link, latest_review_date) shifts one
# position to the left. Deleting the entire column reproduces exactly
# that shift, including updating the used-range dimension from K to J.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E:E").Delete()
